$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "capacity"
$ws.Name = "capacity"

# Update row 2 (Assessment) values
$ws.Range("B2").Value = 11
$ws.Range("C2").Value = 15
$ws.Range("D2").Value = 15
$ws.Range("E2").Value = 35

# Add new service rows
$ws.Range("A3").Value = "MBT"
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = 16
$ws.Range("D3").Value = 21
$ws.Range("E3").Value = 17

$ws.Range("A4").Value = "SCM"
$ws.Range("B4").Value = 12
$ws.Range("C4").Value = 35
$ws.Range("D4").Value = 35
$ws.Range("E4").Value = 20

$ws.Range("A5").Value = "Liaison"
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 8
$ws.Range("D5").Value = 10
$ws.Range("E5").Value = 4

$ws.Range("A6").Value = "Med-RV"
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 8

# New "clock-stop" capacity formulas
$ws.Range("I30").Formula = "=365+250"
$ws.Range("I31").Formula = "=365+30+30+250"

# Reflect the new working selection in the UI
$excel.Goto($ws.Range("A1:E6"))
